# Update imputed values produced by the RandomForest algorithm run
# (commit message: "Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.405400000000004
$ws.Range("B21").Value = 9.518600000000005
$ws.Range("B23").Value = 9.008400000000005
$ws.Range("B25").Value = 6.360699999999996
$ws.Range("C27").Value = -12.7968
$ws.Range("C31").Value = -13.20289999999999
$ws.Range("C39").Value = -12.4154
$ws.Range("C48").Value = -11.33549999999999
$ws.Range("C51").Value = -11.299
$ws.Range("C52").Value = -10.97709999999999
$ws.Range("B53").Value = 5.354799999999999
$ws.Range("C55").Value = -13.69609999999999
$ws.Range("C56").Value = -12.7741
$ws.Range("B57").Value = 4.890799999999997
$ws.Range("C57").Value = -13.73779999999999
$ws.Range("B59").Value = 4.980399999999999
$ws.Range("B69").Value = 5.486599999999992
$ws.Range("C73").Value = -12.37360000000001
$ws.Range("B79").Value = 9.388100000000003
$ws.Range("B83").Value = 5.810699999999998
$ws.Range("C89").Value = -10.5942
$ws.Range("C90").Value = -12.14940000000001
$ws.Range("B93").Value = 5.466000000000001
